{"js": "const replacements = [\n  [\"87\u00d762=5394\", \"31\u00d772=2232\"],\n  [\"25\u00d771=1775\", \"96\u00d717=1632\"],\n  [\"65\u00d784=5460\", \"32\u00d796=3072\"],\n  [\"26\u00d738=988\", \"64\u00d736=2304\"],\n  [\"69\u00d757=3933\", \"69\u00d780=5520\"],\n  [\"31\u00d756=1736\", \"74\u00d743=3182\"],\n  [\"39\u00d744=1716\", \"51\u00d722=1122\"],\n  [\"40\u00d797=3880\", \"17\u00d714=238\"],\n  [\"67\u00d758=3886\", \"65\u00d799=6435\"],\n  [\"24\u00d775=1800\", \"44\u00d751=2244\"],\n  [\"19\u00d733=627\", \"35\u00d786=3010\"],\n  [\"65\u00d785=5525\", \"14\u00d736=504\"],\n  [\"30\u00d771=2130\", \"69\u00d781=5589\"],\n  [\"86\u00d790=7740\", \"34\u00d798=3332\"],\n  [\"90\u00d796=8640\", \"98\u00d760=5880\"],\n  [\"34\u00d736=1224\", \"85\u00d795=8075\"],\n  [\"50\u00d761=3050\", \"21\u00d792=1932\"],\n  [\"21\u00d798=2058\", \"45\u00d723=1035\"],\n  [\"73\u00d789=6497\", \"19\u00d717=323\"],\n  [\"15\u00d723=345\", \"50\u00d741=2050\"],\n  [\"67\u00d752=3484\", \"43\u00d799=4257\"],\n  [\"42\u00d792=3864\", \"95\u00d744=4180\"],\n  [\"75\u00d738=2850\", \"63\u00d742=2646\"],\n  [\"56\u00d795=5320\", \"53\u00d730=1590\"],\n  [\"30\u00d792=2760\", \"37\u00d742=1554\"],\n];\n\nfor (const [oldText, newText] of replacements) {\n  const results = context.document.body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n\n  for (const range of results.items) {\n    range.insertText(newText, Word.InsertLocation.replace);\n  }\n  await context.sync();\n}\n", "ps1": "$d = $word.ActiveDocument\n\n$replacements = @(\n    @{ Old = \"87\u00d762=5394\"; New = \"31\u00d772=2232\" },\n    @{ Old = \"25\u00d771=1775\"; New = \"96\u00d717=1632\" },\n    @{ Old = \"65\u00d784=5460\"; New = \"32\u00d796=3072\" },\n    @{ Old = \"26\u00d738=988\";  New = \"64\u00d736=2304\" },\n    @{ Old = \"69\u00d757=3933\"; New = \"69\u00d780=5520\" },\n    @{ Old = \"31\u00d756=1736\"; New = \"74\u00d743=3182\" },\n    @{ Old = \"39\u00d744=1716\"; New = \"51\u00d722=1122\" },\n    @{ Old = \"40\u00d797=3880\"; New = \"17\u00d714=238\" },\n    @{ Old = \"67\u00d758=3886\"; New = \"65\u00d799=6435\" },\n    @{ Old = \"24\u00d775=1800\"; New = \"44\u00d751=2244\" },\n    @{ Old = \"19\u00d733=627\";  New = \"35\u00d786=3010\" },\n    @{ Old = \"65\u00d785=5525\"; New = \"14\u00d736=504\" },\n    @{ Old = \"30\u00d771=2130\"; New = \"69\u00d781=5589\" },\n    @{ Old = \"86\u00d790=7740\"; New = \"34\u00d798=3332\" },\n    @{ Old = \"90\u00d796=8640\"; New = \"98\u00d760=5880\" },\n    @{ Old = \"34\u00d736=1224\"; New = \"85\u00d795=8075\" },\n    @{ Old = \"50\u00d761=3050\"; New = \"21\u00d792=1932\" },\n    @{ Old = \"21\u00d798=2058\"; New = \"45\u00d723=1035\" },\n    @{ Old = \"73\u00d789=6497\"; New = \"19\u00d717=323\" },\n    @{ Old = \"15\u00d723=345\";  New = \"50\u00d741=2050\" },\n    @{ Old = \"67\u00d752=3484\"; New = \"43\u00d799=4257\" },\n    @{ Old = \"42\u00d792=3864\"; New = \"95\u00d744=4180\" },\n    @{ Old = \"75\u00d738=2850\"; New = \"63\u00d742=2646\" },\n    @{ Old = \"56\u00d795=5320\"; New = \"53\u00d730=1590\" },\n    @{ Old = \"30\u00d792=2760\"; New = \"37\u00d742=1554\" }\n)\n\nforeach ($r in $replacements) {\n    $range = $d.Content\n    $find = $range.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    # Find.Execute(FindText, MatchCase, MatchWholeWord, MatchWildcards,\n    #   MatchSoundsLike, MatchAllWordForms, Forward, Wrap, Format,\n    #   ReplaceWith, Replace)\n    $find.Execute($r.Old, $true, $false, $false, $false, $false, $true, 1, $false, $r.New, 2)\n}\n"}
